$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1930.125
$ws.Range("I32").Value = 1922.25
$ws.Range("K32").Value = 1922.25
$ws.Range("M32").Value = -1596.25

$ws.Range("H41").Value = 966.8333
$ws.Range("J41").Value = 966.8333
$ws.Range("L41").Value = 966.8333
$ws.Range("N41").Value = -1846.8333

$ws.Range("H74").Value = 2976.5
$ws.Range("I74").Value = 2976.5
$ws.Range("K74").Value = 2976.5
$ws.Range("M74").Value = -2040.5

$ws.Range("H77").Value = 2976.5
$ws.Range("I77").Value = 2976.5
$ws.Range("K77").Value = 14882.5
$ws.Range("M77").Value = -10202.5

$ws.Range("H88").Value = 2063.9092
$ws.Range("J88").Value = 2258.3157
$ws.Range("L88").Value = 2258.3157
$ws.Range("N88").Value = -3070.3157

$ws.Range("H91").Value = 2063.9092
$ws.Range("J91").Value = 2258.3157
$ws.Range("L91").Value = 2258.3157
$ws.Range("N91").Value = -5066.3157

$ws.Range("H97").Value = 1000000000
$ws.Range("J97").Value = 1000000000
$ws.Range("L97").Value = 3000000000
$ws.Range("N97").Value = -3000000992

$ws.Range("H132").Value = 5642.2666
$ws.Range("I132").Value = 2438.4443
$ws.Range("K132").Value = 7315.3329
$ws.Range("M132").Value = -4785.3329

$ws.Range("H138").Value = 4199.92
$ws.Range("J138").Value = 4230.05
$ws.Range("L138").Value = 12690.15
$ws.Range("N138").Value = -22970.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2691.6
$ws.Range("I2").Value = 708
$ws.Range("K2").Value = 708
$ws.Range("M2").Value = -595

$ws.Range("H45").Value = 4612.091
$ws.Range("I45").Value = 1746.6
$ws.Range("K45").Value = 1746.6
$ws.Range("M45").Value = -1369.6

$ws.Range("H61").Value = 3767.25
$ws.Range("I61").Value = 3473.3635
$ws.Range("K61").Value = 3473.3635
$ws.Range("M61").Value = -3261.3635

$ws.Range("H74").Value = 882.125
$ws.Range("I74").Value = 882.125
$ws.Range("K74").Value = 882.125
$ws.Range("M74").Value = -8.125

$ws.Range("H77").Value = 882.125
$ws.Range("I77").Value = 882.125
$ws.Range("K77").Value = 4410.625
$ws.Range("M77").Value = -42.625

$ws.Range("H92").Value = 47333.668
$ws.Range("J92").Value = 47333.668
$ws.Range("L92").Value = 47333.668
$ws.Range("N92").Value = -52325.668

$ws.Range("H102").Value = 2390.5833
$ws.Range("I102").Value = 1465.6666
$ws.Range("J102").Value = 5165.3335
$ws.Range("K102").Value = 1465.6666
$ws.Range("L102").Value = 5165.3335
$ws.Range("M102").Value = 156.3334
$ws.Range("N102").Value = -8409.333500000001

$ws.Range("H116").Value = 2691.6
$ws.Range("I116").Value = 708
$ws.Range("K116").Value = 708
$ws.Range("M116").Value = 1586

$ws.Range("H132").Value = 5666
$ws.Range("I132").Value = 5666
$ws.Range("K132").Value = 16998
$ws.Range("M132").Value = -14468

$ws.Range("H136").Value = 3767.25
$ws.Range("I136").Value = 3473.3635
$ws.Range("K136").Value = 10420.0905
$ws.Range("M136").Value = -7870.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2691.6
$ws.Range("I3").Value = 708
$ws.Range("K3").Value = 708
$ws.Range("M3").Value = -594

$ws.Range("H22").Value = 840
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -727

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H86").Value = 4986.625
$ws.Range("I86").Value = 5127.5713
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 5127.5713
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -4004.5713
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 4986.625
$ws.Range("I89").Value = 5127.5713
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 25637.8565
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -20021.8565
$ws.Range("N89").Value = -31232

$ws.Range("H94").Value = 1395.625
$ws.Range("I94").Value = 880.7143
$ws.Range("K94").Value = 880.7143
$ws.Range("M94").Value = -429.7143

$ws.Range("H99").Value = 2578.75
$ws.Range("I99").Value = 2355.625
$ws.Range("K99").Value = 2355.625
$ws.Range("M99").Value = -857.625

$ws.Range("H109").Value = 105995
$ws.Range("J109").Value = 105995
$ws.Range("L109").Value = 105995
$ws.Range("N109").Value = -108769

$ws.Range("H134").Value = 4743.25
$ws.Range("I134").Value = 4489.5
$ws.Range("K134").Value = 13468.5
$ws.Range("M134").Value = -10933.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 107.625
$ws.Range("I19").Value = 107.625
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 107.625
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 62.375
$ws.Range("N19").ClearContents()

$ws.Range("H22").Value = 1302.1666
$ws.Range("I22").Value = 186
$ws.Range("J22").Value = 1403.6364
$ws.Range("K22").Value = 186
$ws.Range("L22").Value = 1403.6364
$ws.Range("M22").Value = 164
$ws.Range("N22").Value = -2103.6364

$ws.Range("H24").Value = 107.625
$ws.Range("I24").Value = 107.625
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 107.625
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 62.375
$ws.Range("N24").ClearContents()

$ws.Range("H58").Value = 7130.778
$ws.Range("I58").Value = 6028.8335
$ws.Range("K58").Value = 6028.8335
$ws.Range("M58").Value = -5825.8335

$ws.Range("H99").Value = 6517.4
$ws.Range("I99").Value = 5431.636
$ws.Range("J99").Value = 9503.25
$ws.Range("K99").Value = 5431.636
$ws.Range("L99").Value = 9503.25
$ws.Range("M99").Value = -3933.636
$ws.Range("N99").Value = -12499.25

$ws.Range("H103").Value = 3833.3333
$ws.Range("I103").Value = 3833.3333
$ws.Range("K103").Value = 3833.3333
$ws.Range("M103").Value = -2661.3333

$ws.Range("H126").Value = 6517.4
$ws.Range("I126").Value = 5431.636
$ws.Range("J126").Value = 9503.25
$ws.Range("K126").Value = 16294.908
$ws.Range("L126").Value = 28509.75
$ws.Range("M126").Value = -13824.908
$ws.Range("N126").Value = -33449.75

$ws.Range("H134").Value = 2460.3809
$ws.Range("I134").Value = 2364.9443
$ws.Range("K134").Value = 7094.8329
$ws.Range("M134").Value = -4559.8329

$ws.Range("H136").Value = 7130.778
$ws.Range("I136").Value = 6028.8335
$ws.Range("K136").Value = 18086.5005
$ws.Range("M136").Value = -15536.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H23").Value = 896.25
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 896.25
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2688.75
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -3158.75

$ws.Range("H119").Value = 514.5
$ws.Range("I119").Value = 514.5
$ws.Range("K119").Value = 1543.5
$ws.Range("M119").Value = 3294.5

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 750
$ws.Range("I9").Value = 750
$ws.Range("K9").Value = 750
$ws.Range("M9").Value = -580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H122").Value = 7399
$ws.Range("I122").Value = 6485.5713
$ws.Range("K122").Value = 19456.7139
$ws.Range("M122").Value = -17006.7139

$ws.Range("H136").Value = 4778.4
$ws.Range("I136").Value = 3997.5
$ws.Range("K136").Value = 11992.5
$ws.Range("M136").Value = -9442.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H81").Value = 534.6
$ws.Range("I81").Value = 728
$ws.Range("J81").Value = 244.5
$ws.Range("K81").Value = 1456
$ws.Range("L81").Value = 489
$ws.Range("M81").Value = -395
$ws.Range("N81").Value = -2611

$ws.Range("H84").Value = 534.6
$ws.Range("I84").Value = 728
$ws.Range("J84").Value = 244.5
$ws.Range("K84").Value = 7280
$ws.Range("L84").Value = 2445
$ws.Range("M84").Value = -1976
$ws.Range("N84").Value = -13053

$ws.Range("H112").Value = 50382.332
$ws.Range("J112").Value = 50382.332
$ws.Range("L112").Value = 50382.332
$ws.Range("N112").Value = -53336.332

$ws.Range("H126").Value = 1561.25
$ws.Range("I126").Value = 1561.25
$ws.Range("K126").Value = 4683.75
$ws.Range("M126").Value = -2213.75

$ws.Range("H136").Value = 1622.7368
$ws.Range("I136").Value = 1102.3871
$ws.Range("J136").Value = 3927.1428
$ws.Range("K136").Value = 3307.1613
$ws.Range("L136").Value = 11781.4284
$ws.Range("M136").Value = -757.1612999999998
$ws.Range("N136").Value = -16881.4284
